# New Testcases: Div Address, Terms Code and Scrap Reason Code
#
# Adds 4 new columns (Backflush Location, Inspection Order Location,
# Vendor Location, Issue Sequence for Backflush) to the
# "Create_Inventory Loc ID" sheet, fills in values for the existing
# row plus two brand-new test-case rows, and tweaks column widths /
# selections on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Create_Inventory Loc ID")
$ws2 = $wb.Worksheets.Item("Edit_Inventory Loc ID")

# --- Sheet1 ("Create_Inventory Loc ID"): new header columns G:J ---
$ws1.Range("G1").Value = "Backflush Location"
$ws1.Range("H1").Value = "Inspection Order Location"
$ws1.Range("I1").Value = "Vendor Location"
$ws1.Range("J1").Value = "Issue Sequence for Backflush"

# --- Row 2: fill in new columns for the existing PK-O test case ---
$ws1.Range("G2").Value = $false
$ws1.Range("H2").Value = $false
$ws1.Range("I2").Value = $false
$ws1.Range("J2").Value = 1

# --- Row 3: new test case (PK-N / Non-nettable) ---
$ws1.Range("A3").Value = "10 (Denver)"
$ws1.Range("B3").Value = "PK-N"
$ws1.Range("C3").Value = "Non-nettable"
$ws1.Range("D3").Value = "PK NN Loc ID"
$ws1.Range("E3").Value = $false
$ws1.Range("F3").Value = $false
$ws1.Range("G3").Value = $false
$ws1.Range("H3").Value = $true
$ws1.Range("I3").Value = $true
$ws1.Range("J3").Value = 0

# --- Row 4: new test case (PK-C / Consigned) ---
$ws1.Range("A4").Value = "10 (Denver)"
$ws1.Range("C4").Value = "Consigned"
$ws1.Range("D4").Value = "PK Cn Loc ID"
$ws1.Range("B4").Value = "PK-C"
$ws1.Range("E4").Value = $true
$ws1.Range("F4").Value = $true
$ws1.Range("G4").Value = $false
$ws1.Range("H4").Value = $false
$ws1.Range("I4").Value = $false
$ws1.Range("J4").Value = 1

# --- Column widths for the new columns (best-fit like the originals) ---
$ws1.Columns.Item(7).ColumnWidth = 15.666666666666666
$ws1.Columns.Item(8).ColumnWidth = 21.666666666666668
$ws1.Columns.Item(9).ColumnWidth = 13.666666666666666
$ws1.Columns.Item(10).ColumnWidth = 24.0

# --- Sheet2 ("Edit_Inventory Loc ID"): widen column B, move selection ---
$ws2.Columns.Item(2).ColumnWidth = 20.0
$ws2.Range("B6").Select()

# --- Leave the workbook focused back on sheet1, scrolled/selected at I3 ---
$ws1.Application.ActiveWindow.ScrollRow = 1
$ws1.Application.ActiveWindow.ScrollColumn = 4
$ws1.Range("I3").Select()
